$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column M by duplicating column L (formatting + values),
# shifting nothing to the right since M was previously empty/out of range.
$ws.Columns("L").Copy()
$ws.Columns("M").Insert(-4161)
$excel.CutCopyMode = $false

# Fix up the values in the new column M for the 2021 data point.
$ws.Range("M4").Value = 2021
$ws.Range("M5").Value = 98
$ws.Range("M6").Value = 97
$ws.Range("M7").Value = 96
